# 1DES FPOO anexado link exemplo SORT
# Fill in column K (the "L1-Log" attendance column) for rows 3-20 on the
# FREQ sheet with the same P/F (Presente/Falta) values recorded in the
# source diff, then leave the active selection on K8 (matching the
# author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FREQ")

$values = @{
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "F"
    8  = "F"
    9  = "P"
    10 = "F"
    11 = "P"
    12 = "P"
    13 = "P"
    14 = "F"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
}

foreach ($row in $values.Keys) {
    $ws.Range("K$row").Value = $values[$row]
}

$ws.Range("K8").Select() | Out-Null
